$wb = $excel.ActiveWorkbook

# Worksheet handles
$wsMain    = $wb.Worksheets.Item("PO_main")
$wsReduced = $wb.Worksheets.Item("PO_reduced")
$ws0to1    = $wb.Worksheets.Item("NPO_0to1")
$ws1to2    = $wb.Worksheets.Item("NPO_1to2")
$ws2to3    = $wb.Worksheets.Item("NPO_2to3")

# ---------------------------------------------------------------------------
# NPO_0to1: strip the " 0|1" interval qualifier from the comparison labels
# ---------------------------------------------------------------------------
$ws0to1.Range("A2").Value  = "LLS vs LUS"
$ws0to1.Range("A3").Value  = "RML vs LUS "
$ws0to1.Range("A4").Value  = "RUL vs LUS"
$ws0to1.Range("A5").Value  = "RLL vs LUS"
$ws0to1.Range("A6").Value  = "LLS vs LLL"
$ws0to1.Range("A7").Value  = "RML vs LLL"
$ws0to1.Range("A8").Value  = "RUL vs LLL "
$ws0to1.Range("A9").Value  = "LLL vs LUS"
$ws0to1.Range("A10").Value = "LLS vs RLL"
$ws0to1.Range("A11").Value = " RML vs RLL"
$ws0to1.Range("A12").Value = "RLL vs LLL"
$ws0to1.Range("A13").Value = "LLS vs RUL"
$ws0to1.Range("A14").Value = "RML vs RUL"
$ws0to1.Range("A15").Value = "RUL vs RLL"
$ws0to1.Range("A16").Value = "LLS vs RML"

# ---------------------------------------------------------------------------
# NPO_1to2: strip the " 1|2" interval qualifier from the comparison labels,
# and clear the stray "applyFill" style that had been left on column A
# ---------------------------------------------------------------------------
$ws1to2.Range("A2").Value  = "RML vs LUS"
$ws1to2.Range("A3").Value  = "LLS vs LUS"
$ws1to2.Range("A4").Value  = "RML vs LLL"
$ws1to2.Range("A5").Value  = "RUL vs LUS"
$ws1to2.Range("A6").Value  = "RLL vs LUS"
$ws1to2.Range("A7").Value  = "RML vs RLL"
$ws1to2.Range("A8").Value  = "RML vs RUL"
$ws1to2.Range("A9").Value  = "LLL vs LUS"
$ws1to2.Range("A10").Value = "RML vs LLS"
$ws1to2.Range("A11").Value = "LLS vs LLL"
$ws1to2.Range("A12").Value = "RUL vs LLL"
$ws1to2.Range("A13").Value = "LLS vs RLL"
$ws1to2.Range("A14").Value = "RLL vs LLL"
$ws1to2.Range("A15").Value = "LLS vs RUL"
$ws1to2.Range("A16").Value = "RUL vs RLL"
$ws1to2.Range("A2:A16").Style = "Normal"

# ---------------------------------------------------------------------------
# NPO_2to3: strip the " 2|3" interval qualifier from the comparison labels,
# and clear the stray "applyFill" style that had been left on column A
# ---------------------------------------------------------------------------
$ws2to3.Range("A2").Value  = "RML vs LUS"
$ws2to3.Range("A3").Value  = "RML vs LLL"
$ws2to3.Range("A4").Value  = "RML vs RLL"
$ws2to3.Range("A5").Value  = "RML vs RUL"
$ws2to3.Range("A6").Value  = "LLS vs LUS"
$ws2to3.Range("A7").Value  = "RUL vs LUS"
$ws2to3.Range("A8").Value  = "RML vs LLS"
$ws2to3.Range("A9").Value  = "RLL vs LUS"
$ws2to3.Range("A10").Value = "LLL vs LUS"
$ws2to3.Range("A11").Value = "LLS vs LLL"
$ws2to3.Range("A12").Value = "RUL vs LLL"
$ws2to3.Range("A13").Value = "LLS vs RLL"
$ws2to3.Range("A14").Value = "RUL vs RLL"
$ws2to3.Range("A15").Value = "RLL vs LLL"
$ws2to3.Range("A16").Value = "LLS vs RUL"
$ws2to3.Range("A2:A16").Style = "Normal"

# ---------------------------------------------------------------------------
# PO_reduced now shares PO_main's shorter label set, so column A is
# auto-fitted down to PO_main's width (the reviewer re-sized it while going
# over the reduced-model tab).
# ---------------------------------------------------------------------------
$wsReduced.Range("A1:A16").EntireColumn.AutoFit()

# ---------------------------------------------------------------------------
# Restore the cursor/selection on every tab, then finish on PO_reduced so it
# becomes the active tab (matches activeTab moving from NPO_2to3 to
# PO_reduced, and tabSelected moving off NPO_2to3 onto PO_reduced).
# ---------------------------------------------------------------------------
$wsMain.Range("A15").Select()
$ws0to1.Range("A2").Select()
$ws1to2.Range("A22").Select()
$ws2to3.Range("A18").Select()
$wsReduced.Activate()
$wsReduced.Range("E19").Select()
